$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (shifts old B->D, old C->E)
$ws.Columns("B:C").Insert()

# New header values for the inserted columns
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new columns with "UN" for every data row (2-27), matching existing B/D/E pattern
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Match column widths (8.0 characters, same as the original column C) for columns C, D and E
$ws.Columns("C").ColumnWidth = 7.1666666666666667
$ws.Columns("D").ColumnWidth = 7.1666666666666667
$ws.Columns("E").ColumnWidth = 7.1666666666666667
